$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 becomes what row 3 used to be, and vice versa (swap D, J, K, L, M, P)
$ws.Range("D2").Value = 44203
$ws.Range("J2").Value = 27
$ws.Range("K2").Value = 7000
$ws.Range("L2").Value = 8000
$ws.Range("M2").Value = 7556
$ws.Range("P2").Value = 756

$ws.Range("D3").Value = 44211
$ws.Range("J3").Value = 28
$ws.Range("K3").Value = 8000
$ws.Range("L3").Value = 8500
$ws.Range("M3").Value = 8214
$ws.Range("P3").Value = 821
